$wb = $excel.ActiveWorkbook

# ALC row 7 (@@ -984,25 +984,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 10975.833
$ws.Range("I7").Value = 6186.6665
$ws.Range("J7").Value = 15765
$ws.Range("K7").Value = 6186.6665
$ws.Range("L7").Value = 15765
$ws.Range("M7").Value = -6074.6665
$ws.Range("N7").Value = -15989

# ALC row 14 (@@ -1333,25 +1333,25 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H14").Value = 10975.833
$ws.Range("I14").Value = 6186.6665
$ws.Range("J14").Value = 15765
$ws.Range("K14").Value = 6186.6665
$ws.Range("L14").Value = 15765
$ws.Range("M14").Value = -5995.6665
$ws.Range("N14").Value = -16147

# ALC row 97 (@@ -5436,22 +5436,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H97").Value = 2702.6365
$ws.Range("J97").Value = 2702.6365
$ws.Range("L97").Value = 8107.9095
$ws.Range("N97").Value = -9099.9095

# ALC row 134 (@@ -7267,22 +7267,22 @@)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H134").Value = 159999
$ws.Range("J134").Value = 159999
$ws.Range("L134").Value = 159999
$ws.Range("N134").Value = -170139

# ARM row 13 (@@ -8301,22 +8301,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H13").Value = 449
$ws.Range("I13").Value = 449
$ws.Range("K13").Value = 449
$ws.Range("M13").Value = -305

# ARM row 17 (@@ -8503,25 +8503,22 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H17").Value = 2125
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 2125
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2125
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -2471

# ARM row 50 (@@ -10153,25 +10150,25 @@)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 834.75
$ws.Range("J50").Value = 90
$ws.Range("L50").Value = 90
$ws.Range("N50").Value = -1518

# BSM row 5 (@@ -14827,25 +14824,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 3400.6
$ws.Range("I5").Value = 2251.5
$ws.Range("J5").Value = 4166.6665
$ws.Range("K5").Value = 2251.5
$ws.Range("L5").Value = 4166.6665
$ws.Range("M5").Value = -2138.5
$ws.Range("N5").Value = -4392.6665

# BSM row 11 (@@ -15121,25 +15118,25 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 3750
$ws.Range("J11").Value = 5000
$ws.Range("L11").Value = 5000
$ws.Range("N11").Value = -5280

# BSM row 12 (@@ -15173,25 +15170,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H12").Value = 4999
$ws.Range("I12").Value = 4999
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 4999
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = -4831
$ws.Range("N12").ClearContents()

# BSM row 38 (@@ -16432,22 +16426,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 39789
$ws.Range("J38").Value = 39789
$ws.Range("L38").Value = 39789
$ws.Range("N38").Value = -40621

# BSM row 44 (@@ -16717,22 +16711,22 @@)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value = 48662.332
$ws.Range("J44").Value = 48662.332
$ws.Range("L44").Value = 48662.332
$ws.Range("N44").Value = -49656.332

# CRP row 2 (@@ -21547,25 +21541,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 467.625
$ws.Range("J2").Value = 651.4
$ws.Range("L2").Value = 651.4
$ws.Range("N2").Value = -877.4

# CRP row 4 (@@ -21648,25 +21642,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 74367336
$ws.Range("I4").Value = 32999
$ws.Range("J4").Value = 353121100
$ws.Range("K4").Value = 32999
$ws.Range("L4").Value = 353121100
$ws.Range("M4").Value = -32887
$ws.Range("N4").Value = -353121324

# CRP row 6 (@@ -21752,26 +21746,23 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 348666
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# CRP row 19 (@@ -22401,25 +22392,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 1898
$ws.Range("I19").Value = 973.5
$ws.Range("J19").Value = 2324.6924
$ws.Range("K19").Value = 973.5
$ws.Range("L19").Value = 2324.6924
$ws.Range("M19").Value = -803.5
$ws.Range("N19").Value = -2664.6924

# CRP row 24 (@@ -22649,25 +22640,25 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H24").Value = 1898
$ws.Range("I24").Value = 973.5
$ws.Range("J24").Value = 2324.6924
$ws.Range("K24").Value = 973.5
$ws.Range("L24").Value = 2324.6924
$ws.Range("M24").Value = -803.5
$ws.Range("N24").Value = -2664.6924

# CRP row 25 (@@ -22701,22 +22692,22 @@)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 6153.4287
$ws.Range("I25").Value = 1007.3333
$ws.Range("K25").Value = 1007.3333
$ws.Range("M25").Value = -833.3333

# CUL row 18 (@@ -29327,25 +29318,25 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 223.75
$ws.Range("J18").Value = 133
$ws.Range("L18").Value = 399
$ws.Range("N18").Value = -737

# CUL row 140 (@@ -35416,22 +35407,22 @@)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 4681.6113
$ws.Range("I140").Value = 1847.6364
$ws.Range("K140").Value = 5542.9092
$ws.Range("M140").Value = -362.9092000000001

# GSM row 6 (@@ -35816,19 +35807,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("L6").Value = 1
$ws.Range("N6").Value = -227

# GSM row 11 (@@ -36049,25 +36043,25 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 1117284
$ws.Range("I11").Value = 1000001
$ws.Range("J11").Value = 1234567
$ws.Range("K11").Value = 1000001
$ws.Range("L11").Value = 1234567
$ws.Range("M11").Value = -999862
$ws.Range("N11").Value = -1234845

# GSM row 16 (@@ -36291,19 +36285,22 @@)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H16").Value = 1
$ws.Range("J16").Value = 1
$ws.Range("L16").Value = 1
$ws.Range("N16").Value = -501

# LTW row 13 (@@ -43032,19 +43029,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H13").Value = 3000
$ws.Range("J13").Value = 3000
$ws.Range("L13").Value = 3000
$ws.Range("N13").Value = -3280

# LTW row 19 (@@ -43329,25 +43329,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 4823.778
$ws.Range("I19").Value = 2003.75
$ws.Range("J19").Value = 7079.8
$ws.Range("K19").Value = 2003.75
$ws.Range("L19").Value = 7079.8
$ws.Range("M19").Value = -1833.75
$ws.Range("N19").Value = -7419.8

# LTW row 20 (@@ -43381,25 +43381,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 84113.06
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 84113.06
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 84113.06
$ws.Range("M20").ClearContents()
$ws.Range("N20").Value = -84565.06

# LTW row 51 (@@ -44933,22 +44930,22 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H51").Value = 34143
$ws.Range("J51").Value = 34143
$ws.Range("L51").Value = 34143
$ws.Range("N51").Value = -35099

# LTW row 53 (@@ -45031,25 +45028,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H53").Value = 6114
$ws.Range("J53").Value = 5724.5
$ws.Range("L53").Value = 5724.5
$ws.Range("N53").Value = -6760.5

# LTW row 136 (@@ -49074,25 +49071,25 @@)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 4176.615
$ws.Range("I136").Value = 3890.6365
$ws.Range("J136").Value = 5749.5
$ws.Range("K136").Value = 11671.9095
$ws.Range("L136").Value = 17248.5
$ws.Range("M136").Value = -9121.9095
$ws.Range("N136").Value = -22348.5

# WVR row 4 (@@ -49554,22 +49551,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 949.5
$ws.Range("I4").Value = 949.3333
$ws.Range("K4").Value = 949.3333
$ws.Range("M4").Value = -836.3333

# WVR row 8 (@@ -49753,25 +49750,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H8").Value = 2000
$ws.Range("J8").Value = 3000
$ws.Range("L8").Value = 3000
$ws.Range("N8").Value = -3280

# WVR row 11 (@@ -49900,23 +49897,20 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 5000
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

# WVR row 17 (@@ -50191,22 +50185,22 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 6166
$ws.Range("I17").Value = 6166
$ws.Range("K17").Value = 6166
$ws.Range("M17").Value = -5994

# WVR row 132 (@@ -55820,25 +55814,25 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 19142.443
$ws.Range("I132").Value = 11527.477
$ws.Range("J132").Value = 35975.527
$ws.Range("K132").Value = 34582.431
$ws.Range("L132").Value = 107926.581
$ws.Range("M132").Value = -32052.431
$ws.Range("N132").Value = -112986.581

# WVR row 133 (@@ -55872,22 +55866,19 @@)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()
